$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.57568926279046
$ws.Range("C2").Value = 8.411402818432384
$ws.Range("D2").Value = 5.77987598320147
$ws.Range("E2").Value = 12.338066472155
$ws.Range("F2").Value = 27.5793829702987
$ws.Range("I2").Value = 24.88168758446635
$ws.Range("K2").Value = 8.990456234174768
$ws.Range("L2").Value = 10.1392141833351
$ws.Range("M2").Value = 13.85074040336107
$ws.Range("N2").Value = 20.01073722370116
$ws.Range("O2").Value = 24.7763701348219

$ws.Range("B3").Value = 11.35274404535038
$ws.Range("C3").Value = 8.368775653057654
$ws.Range("D3").Value = 5.726312211540399
$ws.Range("E3").Value = 12.3646028003379
$ws.Range("F3").Value = 27.6151517684014
$ws.Range("I3").Value = 24.96151525016285
$ws.Range("K3").Value = 8.824049441003243
$ws.Range("L3").Value = 10.14659296349581
$ws.Range("M3").Value = 13.81954460307021
$ws.Range("N3").Value = 20.06815077141586
$ws.Range("O3").Value = 24.84301352394798

$ws.Range("B4").Value = 11.2156673369071
$ws.Range("C4").Value = 8.34211231900327
$ws.Range("D4").Value = 5.692592784133676
$ws.Range("E4").Value = 12.38243338315651
$ws.Range("F4").Value = 27.64333794687496
$ws.Range("I4").Value = 25.01454087973668
$ws.Range("K4").Value = 8.721702006289339
$ws.Range("L4").Value = 10.1524957887733
$ws.Range("M4").Value = 13.80233469023047
$ws.Range("N4").Value = 20.10506738927893
$ws.Range("O4").Value = 24.88844642166266

$ws.Range("B5").Value = 11.15983759346573
$ws.Range("C5").Value = 8.331126388260882
$ws.Range("D5").Value = 5.678648164320137
$ws.Range("E5").Value = 12.39008638610733
$ws.Range("F5").Value = 27.65638779086868
$ws.Range("I5").Value = 25.03715771493138
$ws.Range("K5").Value = 8.680006464146574
$ws.Range("L5").Value = 10.1552468966056
$ws.Range("M5").Value = 13.79581564970971
$ws.Range("N5").Value = 20.12053095236793
$ws.Range("O5").Value = 24.90809458980674

$ws.Range("B6").Value = 11.15057121412935
$ws.Range("C6").Value = 8.329294993222891
$ws.Range("D6").Value = 5.676320525451258
$ws.Range("E6").Value = 12.39138053936738
$ws.Range("F6").Value = 27.65864910460463
$ws.Range("I6").Value = 25.04097412678637
$ws.Range("K6").Value = 8.673085305470837
$ws.Range("L6").Value = 10.15572460911252
$ws.Range("M6").Value = 13.79476315752902
$ws.Range("N6").Value = 20.12312405184793
$ws.Range("O6").Value = 24.91142560415388

$ws.Range("B7").Value = 11.21491417227249
$ws.Range("C7").Value = 8.341964642681015
$ws.Range("D7").Value = 5.692405539756534
$ws.Range("E7").Value = 12.38253502722681
$ws.Range("F7").Value = 27.64350761223989
$ws.Range("I7").Value = 25.01484181529231
$ws.Range("K7").Value = 8.721139563680213
$ws.Range("L7").Value = 10.15253149087051
$ws.Range("M7").Value = 13.80224476482393
$ws.Range("N7").Value = 20.10527423519203
$ws.Range("O7").Value = 24.88870681413945

$ws.Range("B8").Value = 11.49890249683494
$ws.Range("C8").Value = 8.39680838124289
$ws.Range("D8").Value = 5.76158282640305
$ws.Range("E8").Value = 12.34689737920939
$ws.Range("F8").Value = 27.59042362197627
$ws.Range("I8").Value = 24.90837955863764
$ws.Range("K8").Value = 8.933147954733627
$ws.Range("L8").Value = 10.14147396146906
$ws.Range("M8").Value = 13.83958374353299
$ws.Range("N8").Value = 20.03018867903767
$ws.Range("O8").Value = 24.79841135351477

$ws.Range("B9").Value = 12.05101497620471
$ws.Range("C9").Value = 8.500346187602585
$ws.Range("D9").Value = 5.890419064262781
$ws.Range("E9").Value = 12.28919261571202
$ws.Range("F9").Value = 27.53574918551134
$ws.Range("I9").Value = 24.73144106361075
$ws.Range("K9").Value = 9.345153287726948
$ws.Range("L9").Value = 10.13064981085861
$ws.Range("M9").Value = 13.9279913284089
$ws.Range("N9").Value = 19.89609731982589
$ws.Range("O9").Value = 24.65720150289095

$ws.Range("B10").Value = 12.44948469270492
$ws.Range("C10").Value = 8.573810382262257
$ws.Range("D10").Value = 5.980624122030444
$ws.Range("E10").Value = 12.25420025585162
$ws.Range("F10").Value = 27.52573124378134
$ws.Range("I10").Value = 24.62086040735545
$ws.Range("K10").Value = 9.642513765553016
$ws.Range("L10").Value = 10.12927665151721
$ws.Range("M10").Value = 14.0018574982666
$ws.Range("N10").Value = 19.8055185010306
$ws.Range("O10").Value = 24.57537437890571

$ws.Range("B11").Value = 12.62836285587491
$ws.Range("C11").Value = 8.606631694091606
$ws.Range("D11").Value = 6.020634394715334
$ws.Range("E11").Value = 12.23988393843665
$ws.Range("F11").Value = 27.52771327147661
$ws.Range("I11").Value = 24.5747724245505
$ws.Range("K11").Value = 9.776024249251549
$ws.Range("L11").Value = 10.13007038775864
$ws.Range("M11").Value = 14.03731621514009
$ws.Range("N11").Value = 19.76601863288473
$ws.Range("O11").Value = 24.54291948694569

$ws.Range("B12").Value = 12.69569073700109
$ws.Range("C12").Value = 8.618971370480468
$ws.Range("D12").Value = 6.035632995890083
$ws.Range("E12").Value = 12.23469265243313
$ws.Range("F12").Value = 27.52940236268622
$ws.Range("I12").Value = 24.5579265447453
$ws.Range("K12").Value = 9.826280686480748
$ws.Range("L12").Value = 10.13057390198598
$ws.Range("M12").Value = 14.05100315682212
$ws.Range("N12").Value = 19.75130497737673
$ws.Range("O12").Value = 24.53131610064592

$ws.Range("B13").Value = 12.68120969209698
$ws.Range("C13").Value = 8.616317809114515
$ws.Range("D13").Value = 6.032409637484601
$ws.Range("E13").Value = 12.2358004656101
$ws.Range("F13").Value = 27.52899688040979
$ws.Range("I13").Value = 24.56152761930313
$ws.Range("K13").Value = 9.815471189816428
$ws.Range("L13").Value = 10.13045645233561
$ws.Range("M13").Value = 14.04804401029008
$ws.Range("N13").Value = 19.7544629874927
$ws.Range("O13").Value = 24.53378454962532

$ws.Range("B14").Value = 12.63391053081448
$ws.Range("C14").Value = 8.607648688704613
$ws.Range("D14").Value = 6.021871422133864
$ws.Range("E14").Value = 12.23945224091158
$ws.Range("F14").Value = 27.5278334368958
$ws.Range("I14").Value = 24.57337434353995
$ws.Range("K14").Value = 9.780165176083853
$ws.Range("L14").Value = 10.13010775102006
$ws.Range("M14").Value = 14.03843709133963
$ws.Range("N14").Value = 19.76480324583234
$ws.Range("O14").Value = 24.54195110493225

$ws.Range("B15").Value = 12.60488321821594
$ws.Range("C15").Value = 8.602326915662857
$ws.Range("D15").Value = 6.015396462184384
$ws.Range("E15").Value = 12.24171899973136
$ws.Range("F15").Value = 27.52724295533508
$ws.Range("I15").Value = 24.5807098189176
$ws.Range("K15").Value = 9.758498634766916
$ws.Range("L15").Value = 10.12992055896487
$ws.Range("M15").Value = 14.03258614273715
$ws.Range("N15").Value = 19.77116870831139
$ws.Range("O15").Value = 24.54704279149088

$ws.Range("B16").Value = 12.43774080459659
$ws.Range("C16").Value = 8.57165313853451
$ws.Range("D16").Value = 5.977988302290825
$ws.Range("E16").Value = 12.25516806475391
$ws.Range("F16").Value = 27.52573316143253
$ws.Range("I16").Value = 24.6239572329016
$ws.Range("K16").Value = 9.633748960458531
$ws.Range("L16").Value = 10.12925323540449
$ws.Range("M16").Value = 13.99957691484022
$ws.Range("N16").Value = 19.80813412502882
$ws.Range("O16").Value = 24.57759141557742

$ws.Range("B17").Value = 12.33454460531688
$ws.Range("C17").Value = 8.552680444375838
$ws.Range("D17").Value = 5.954773542777229
$ws.Range("E17").Value = 12.26382866264432
$ws.Range("F17").Value = 27.52648073939821
$ws.Range("I17").Value = 24.651568171473
$ws.Range("K17").Value = 9.556733273345991
$ws.Range("L17").Value = 10.12920643922653
$ws.Range("M17").Value = 13.97979712026143
$ws.Range("N17").Value = 19.83124710687382
$ws.Range("O17").Value = 24.59755402069108

$ws.Range("B18").Value = 12.27496785318847
$ws.Range("C18").Value = 8.541711755757731
$ws.Range("D18").Value = 5.94132489667486
$ws.Range("E18").Value = 12.26896080853233
$ws.Range("F18").Value = 27.52752637944086
$ws.Range("I18").Value = 24.66784603111694
$ws.Range("K18").Value = 9.512272906571283
$ws.Range("L18").Value = 10.12931309541837
$ws.Range("M18").Value = 13.96859552769376
$ws.Range("N18").Value = 19.84470161863103
$ws.Range("O18").Value = 24.60948479051295

$ws.Range("B19").Value = 12.25476030140921
$ws.Range("C19").Value = 8.537988405168662
$ws.Range("D19").Value = 5.936755063520695
$ws.Range("E19").Value = 12.27072437542125
$ws.Range("F19").Value = 27.52798619812367
$ws.Range("I19").Value = 24.67342556859747
$ws.Range("K19").Value = 9.497192902396176
$ws.Range("L19").Value = 10.12937217822944
$ws.Range("M19").Value = 13.96483317857593
$ws.Range("N19").Value = 19.84928469023228
$ws.Range("O19").Value = 24.61360139948956

$ws.Range("B20").Value = 12.34555338294724
$ws.Range("C20").Value = 8.554705948100786
$ws.Range("D20").Value = 5.957254780113272
$ws.Range("E20").Value = 12.26289112261357
$ws.Range("F20").Value = 27.5263374515849
$ws.Range("I20").Value = 24.64858787176741
$ws.Range("K20").Value = 9.564948948323883
$ws.Range("L20").Value = 10.12919760367435
$ws.Range("M20").Value = 13.9818846315917
$ws.Range("N20").Value = 19.82877008344151
$ws.Range("O20").Value = 24.59538251089247

$ws.Range("B21").Value = 12.64781505934501
$ws.Range("C21").Value = 8.610197460400224
$ws.Range("D21").Value = 6.024970928064861
$ws.Range("E21").Value = 12.23837338696547
$ws.Range("F21").Value = 27.5281497143026
$ws.Range("I21").Value = 24.56987820612726
$ws.Range("K21").Value = 9.790543948604496
$ws.Range("L21").Value = 10.13020467375555
$ws.Range("M21").Value = 14.04125189482573
$ws.Range("N21").Value = 19.76175944535382
$ws.Range("O21").Value = 24.53953375003616

$ws.Range("B22").Value = 12.842946347375
$ws.Range("C22").Value = 8.645943959608477
$ws.Range("D22").Value = 6.068336654258352
$ws.Range("E22").Value = 12.22369002646547
$ws.Range("F22").Value = 27.53480336634702
$ws.Range("I22").Value = 24.5219730345739
$ws.Range("K22").Value = 9.936208007694173
$ws.Range("L22").Value = 10.13204523118396
$ws.Range("M22").Value = 14.08156069007253
$ws.Range("N22").Value = 19.71938635626613
$ws.Range("O22").Value = 24.50703555739868

$ws.Range("B23").Value = 12.73904283088701
$ws.Range("C23").Value = 8.626914004595051
$ws.Range("D23").Value = 6.045274697716629
$ws.Range("E23").Value = 12.23140428758428
$ws.Range("F23").Value = 27.53075251474535
$ws.Range("I23").Value = 24.54721723150655
$ws.Range("K23").Value = 9.858641955212473
$ws.Range("L23").Value = 10.13095506149628
$ws.Range("M23").Value = 14.05991161629618
$ws.Range("N23").Value = 19.74187188874484
$ws.Range("O23").Value = 24.52401400486011

$ws.Range("B24").Value = 12.34057708161128
$ws.Range("C24").Value = 8.553790407755601
$ws.Range("D24").Value = 5.95613333056682
$ws.Range("E24").Value = 12.26331450769868
$ws.Range("F24").Value = 27.5264003135638
$ws.Range("I24").Value = 24.64993400697641
$ws.Range("K24").Value = 9.561235208060246
$ws.Range("L24").Value = 10.12920118213013
$ws.Range("M24").Value = 13.98094033739106
$ws.Range("N24").Value = 19.82988942685221
$ws.Range("O24").Value = 24.59636283634952

$ws.Range("B25").Value = 11.90261629527827
$ws.Range("C25").Value = 8.472782361473625
$ws.Range("D25").Value = 5.856326282880088
$ws.Range("E25").Value = 12.3035015198903
$ws.Range("F25").Value = 27.54524434100062
$ws.Range("I25").Value = 24.77589860881547
$ws.Range("K25").Value = 9.234418912979383
$ws.Range("L25").Value = 10.13241935430597
$ws.Range("M25").Value = 13.90248317618705
$ws.Range("N25").Value = 19.93097288709414
$ws.Range("O25").Value = 24.69155703179003
